$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E11").Value = "ASK TO SIR, "
$ws.Range("E11").Select()
